$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Write a value as text, preserving leading zeros, without altering the
    # cell's number format / style (keeps it at the default style index).
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# 1) Update existing row 9 (008028807 / RAFAEL / 4000) -> (004212438 / KENIA / 3750)
Set-TextCell $ws.Cells.Item(9, 1) "004212438"
Set-TextCell $ws.Cells.Item(9, 2) "KENIA"
$ws.Cells.Item(9, 3).Value = 3750

# 2) Insert a new row before row 12 (004392159 / RODRIGO) for KELMA
$ws.Rows.Item(12).Insert()
Set-TextCell $ws.Cells.Item(12, 1) "004504449"
Set-TextCell $ws.Cells.Item(12, 2) "KELMA"
$ws.Cells.Item(12, 3).Value = 1000

# 3) Insert a new row after row 17 (004935287 / ODILON), i.e. before row 18 (MARINA), for MAGALI
$ws.Rows.Item(18).Insert()
Set-TextCell $ws.Cells.Item(18, 1) "004207641"
Set-TextCell $ws.Cells.Item(18, 2) "MAGALI"
$ws.Cells.Item(18, 3).Value = 250

$excel.CutCopyMode = $false
